# Updates cryptos list: new Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "1.000", "30.337.47") rather
# than numbers, so values must be written as literal text. Force Text
# format on every Price cell being updated before assigning its new
# value -- otherwise Excel parses numeric-looking strings as numbers and
# silently drops meaningful trailing zeros (e.g. "1.000" -> 1).
$priceCells = @(
    "D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12",
    "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22",
    "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32",
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42",
    "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.337.47"
$ws.Range("E2").Value = "  -2.83%  "

$ws.Range("D3").Value = "1.936.45"
$ws.Range("E3").Value = "  -2.98%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "251.04"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("D6").Value = "0.7237"
$ws.Range("E6").Value = "  -7.01%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "0.3312"
$ws.Range("E8").Value = "  -4.96%  "

$ws.Range("D9").Value = "28.04"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").Value = "0.07208"
$ws.Range("E10").Value = "  +2.04%  "

$ws.Range("D11").Value = "0.8117"
$ws.Range("E11").Value = "  -3.47%  "

$ws.Range("D12").Value = "0.08094"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("D13").Value = "1.937.47"
$ws.Range("E13").Value = "  -2.90%  "

$ws.Range("D14").Value = "5.491"
$ws.Range("E14").Value = "  -2.49%  "

$ws.Range("D15").Value = "94.66"
$ws.Range("E15").Value = "  -5.99%  "

$ws.Range("D16").Value = "15.20"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "30.356.47"
$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").Value = "0.000008308"
$ws.Range("E18").Value = "  +3.83%  "

$ws.Range("D19").Value = "250.57"
$ws.Range("E19").Value = "  -7.96%  "

$ws.Range("D20").Value = "5.908"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").Value = "2.193.40"
$ws.Range("E21").Value = "  -2.87%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "6.990"
$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").Value = "9.750"
$ws.Range("E25").Value = "  -2.49%  "

$ws.Range("D26").Value = "163.59"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").Value = "2.379"
$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("D28").Value = "19.31"
$ws.Range("E28").Value = "  -2.86%  "

$ws.Range("D29").Value = "0.1324"
$ws.Range("E29").Value = "  -7.12%  "

$ws.Range("D30").Value = "1.571"
$ws.Range("E30").Value = "  -1.55%  "

$ws.Range("D31").Value = "1.352"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").Value = "4.439"
$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("D33").Value = "4.176"
$ws.Range("E33").Value = "  -5.94%  "

$ws.Range("D34").Value = "0.05202"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").Value = "1.286"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("D36").Value = "0.7509"
$ws.Range("E36").Value = "  -5.24%  "

$ws.Range("D37").Value = "2.747"
$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("D38").Value = "0.01981"
$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("D39").Value = "2.834"
$ws.Range("E39").Value = "  -2.72%  "

$ws.Range("D40").Value = "80.06"
$ws.Range("E40").Value = "  -4.01%  "

$ws.Range("D41").Value = "6.442"
$ws.Range("E41").Value = "  -5.09%  "

$ws.Range("D42").Value = "0.4540"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("D43").Value = "2.032"
$ws.Range("E43").Value = "  -4.90%  "

$ws.Range("D44").Value = "0.8482"
$ws.Range("E44").Value = "  -1.02%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").Value = "101.92"
$ws.Range("E46").Value = "  -3.12%  "

$ws.Range("D47").Value = "9.808"
$ws.Range("E47").Value = "  -2.00%  "

$ws.Range("D48").Value = "7.460"
$ws.Range("E48").Value = "  -3.30%  "

$ws.Range("E49").Value = "  -1.70%  "

$ws.Range("D50").Value = "0.4183"
$ws.Range("E50").Value = "  -3.42%  "

$ws.Range("D51").Value = "0.06045"
$ws.Range("E51").Value = "  +0.51%  "
